$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A13").Value = "Buddy Hield"
$ws.Range("B13").Value = "SG,SF"
$ws.Range("C13").Value = "Golden State Warriors"
